$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to remain text so numeric-looking values are not
# auto-converted to numbers by Excel (the source data is text in the workbook).
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '26.301.54'
$ws.Range("E2").Value = '  -0.48%  '
$ws.Range("D3").Value = '1.670.03'
$ws.Range("E3").Value = '  -0.44%  '
$ws.Range("D4").Value = '1.005'
$ws.Range("E4").Value = '  +0.39%  '
$ws.Range("D5").Value = '219.07'
$ws.Range("E5").Value = '  +1.13%  '
$ws.Range("D6").Value = '0.5254'
$ws.Range("E6").Value = '  -1.10%  '
$ws.Range("D7").Value = '1.005'
$ws.Range("E7").Value = '  +0.35%  '
$ws.Range("D8").Value = '0.2678'
$ws.Range("E8").Value = '  -0.69%  '
$ws.Range("D9").Value = '0.06341'
$ws.Range("E9").Value = '  -0.87%  '
$ws.Range("D10").Value = '21.08'
$ws.Range("E10").Value = '  -2.78%  '
$ws.Range("D11").Value = '0.07771'
$ws.Range("E11").Value = '  -0.50%  '
$ws.Range("D12").Value = '1.670.56'
$ws.Range("E12").Value = '  -0.22%  '
$ws.Range("D13").Value = '4.457'
$ws.Range("E13").Value = '  -1.13%  '
$ws.Range("D14").Value = '1.893.79'
$ws.Range("E14").Value = '  -0.64%  '
$ws.Range("D15").Value = '0.5578'
$ws.Range("E15").Value = '  +0.31%  '
$ws.Range("D16").Value = '0.0₅8289'
$ws.Range("E16").Value = '  -0.36%  '
$ws.Range("D17").Value = '65.18'
$ws.Range("E17").Value = '  -0.59%  '
$ws.Range("D18").Value = '26.316.35'
$ws.Range("E18").Value = '  -0.60%  '
$ws.Range("D20").Value = '4.687'
$ws.Range("E20").Value = '  -0.96%  '
$ws.Range("D21").Value = '195.74'
$ws.Range("E21").Value = '  +1.18%  '
$ws.Range("D22").Value = '10.19'
$ws.Range("E22").Value = '  -0.97%  '
$ws.Range("D23").Value = '6.105'
$ws.Range("E23").Value = '  -3.73%  '
$ws.Range("D24").Value = '1.007'
$ws.Range("E24").Value = '  +0.52%  '
$ws.Range("D25").Value = '140.01'
$ws.Range("E25").Value = '  -1.58%  '
$ws.Range("D26").Value = '0.1244'
$ws.Range("E26").Value = '  -3.45%  '
$ws.Range("D27").Value = '7.236'
$ws.Range("E27").Value = '  -2.16%  '
$ws.Range("D28").Value = '16.25'
$ws.Range("E28").Value = '  +0.19%  '
$ws.Range("D29").Value = '1.415'
$ws.Range("E29").Value = '  -1.48%  '
$ws.Range("D30").Value = '0.06192'
$ws.Range("E30").Value = '  -1.03%  '
$ws.Range("D31").Value = '1.283'
$ws.Range("E31").Value = '  +0.77%  '
$ws.Range("D32").Value = '3.609'
$ws.Range("E32").Value = '  -0.03%  '
$ws.Range("D33").Value = '3.313'
$ws.Range("E33").Value = '  -3.86%  '
$ws.Range("D34").Value = '1.640'
$ws.Range("E34").Value = '  -2.21%  '
$ws.Range("D35").Value = '0.9767'
$ws.Range("E35").Value = '  -3.00%  '
$ws.Range("E36").Value = '  -0.03%  '
$ws.Range("D37").Value = '2.786'
$ws.Range("E37").Value = '  +0.17%  '
$ws.Range("D38").Value = '0.5796'
$ws.Range("E38").Value = '  -5.97%  '
$ws.Range("D39").Value = '0.01611'
$ws.Range("E39").Value = '  -1.05%  '
$ws.Range("D40").Value = '6.048'
$ws.Range("E40").Value = '  -1.61%  '
$ws.Range("D41").Value = '0.8593'
$ws.Range("E41").Value = '  -0.51%  '
$ws.Range("E42").Value = '  +0.43%  '
$ws.Range("D43").Value = '1.026.48'
$ws.Range("E43").Value = '  -5.04%  '
$ws.Range("D44").Value = '100.38'
$ws.Range("E44").Value = '  +0.15%  '
$ws.Range("D45").Value = '1.809.48'
$ws.Range("E45").Value = '  -0.69%  '
$ws.Range("D46").Value = '57.96'
$ws.Range("E46").Value = '  +1.46%  '
$ws.Range("B47").Value = 'BabyDogeCoin'
$ws.Range("C47").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D47").Value = '0.0₈109'
$ws.Range("E47").Value = '  +5.46%  '
$ws.Range("B48").Value = 'Frax'
$ws.Range("C48").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D48").Value = '1.009'
$ws.Range("E48").Value = '  +1.30%  '
$ws.Range("B49").Value = 'EnergySwap'
$ws.Range("C49").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D49").Value = '8.073'
$ws.Range("E49").Value = '  -0.75%  '
$ws.Range("D50").Value = '1.492'
$ws.Range("E50").Value = '  +1.29%  '
$ws.Range("B51").Value = 'Cronos'
$ws.Range("C51").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D51").Value = '0.05186'
